$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: average of the |S*|/n column (J) ---
$j12 = $ws.Range("J12")
$j12.Formula = "=AVERAGE(J2:J11)"
$j12.Font.Bold = $true

# --- Rows 14-17: summary labels + stats ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"

$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$bRange = $ws.Range("B14:B17")
$bRange.Font.Bold = $true
$bRange.Font.Size = 12
$bRange.VerticalAlignment = -4108

# --- Page setup (portrait, small paper) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Leave the selection on the new average cell ---
[void]$ws.Range("J12").Select()
